# Auto-generated Excel COM-interop script updating the cryptos price table
# (values refreshed by the scheduled GitHub Actions job).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.528.13"
$ws.Range("E2").Value = "  +2.11%  "
$ws.Range("D3").Value = "1.683.46"
$ws.Range("E3").Value = "  +2.81%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "217.63"
$ws.Range("E5").Value = "  +4.17%  "
$ws.Range("D6").Value = "0.5321"
$ws.Range("E6").Value = "  +3.22%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("D8").Value = "0.2676"
$ws.Range("E8").Value = "  +4.50%  "
$ws.Range("D9").Value = "0.06424"
$ws.Range("E9").Value = "  +3.47%  "
$ws.Range("D10").Value = "21.48"
$ws.Range("E10").Value = "  +5.67%  "
$ws.Range("D11").Value = "0.07774"
$ws.Range("E11").Value = "  +3.33%  "
$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D12").Value = "4.502"
$ws.Range("E12").Value = "  +3.31%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.671.03"
$ws.Range("E13").Value = "  +2.21%  "
$ws.Range("D14").Value = "0.5619"
$ws.Range("E14").Value = "  +4.48%  "
$ws.Range("D15").Value = "0.0₅8402"
$ws.Range("E15").Value = "  +5.64%  "
$ws.Range("D16").Value = "65.97"
$ws.Range("E16").Value = "  +1.78%  "
$ws.Range("D17").Value = "26.566.12"
$ws.Range("E17").Value = "  +2.34%  "
$ws.Range("E18").Value = "  -0.12%  "
$ws.Range("D19").Value = "4.814"
$ws.Range("E19").Value = "  +3.72%  "
$ws.Range("D20").Value = "194.46"
$ws.Range("E20").Value = "  +4.63%  "
$ws.Range("D21").Value = "10.41"
$ws.Range("E21").Value = "  +4.33%  "
$ws.Range("D22").Value = "6.396"
$ws.Range("E22").Value = "  +5.37%  "
$ws.Range("E23").Value = "  -0.05%  "
$ws.Range("D24").Value = "'143.90"
$ws.Range("E24").Value = "  -1.33%  "
$ws.Range("D25").Value = "0.1269"
$ws.Range("E25").Value = "  +6.32%  "
$ws.Range("D26").Value = "7.476"
$ws.Range("E26").Value = "  +2.01%  "
$ws.Range("D27").Value = "'16.20"
$ws.Range("E27").Value = "  +4.72%  "
$ws.Range("D28").Value = "1.414"
$ws.Range("E28").Value = "  +3.18%  "
$ws.Range("D29").Value = "0.06107"
$ws.Range("E29").Value = "  +2.66%  "
$ws.Range("D30").Value = "1.279"
$ws.Range("E30").Value = "  +3.25%  "
$ws.Range("D31").Value = "3.607"
$ws.Range("E31").Value = "  +7.49%  "
$ws.Range("D32").Value = "3.463"
$ws.Range("E32").Value = "  +3.58%  "
$ws.Range("E33").Value = "  +6.29%  "
$ws.Range("D34").Value = "1.017"
$ws.Range("E34").Value = "  +5.03%  "
$ws.Range("D35").Value = "2.795"
$ws.Range("E35").Value = "  +2.95%  "
$ws.Range("D36").Value = "'2.420"
$ws.Range("E36").Value = "  +1.82%  "
$ws.Range("D37").Value = "0.5711"
$ws.Range("E37").Value = "  -1.63%  "
$ws.Range("D38").Value = "0.01639"
$ws.Range("E38").Value = "  +3.77%  "
$ws.Range("D39").Value = "5.969"
$ws.Range("E39").Value = "  +3.45%  "
$ws.Range("D40").Value = "0.8668"
$ws.Range("E40").Value = "  +3.42%  "
$ws.Range("D41").Value = "1.056.30"
$ws.Range("E41").Value = "  +0.12%  "
$ws.Range("E42").Value = "  +0.10%  "
$ws.Range("D43").Value = "100.17"
$ws.Range("E43").Value = "  +0.42%  "
$ws.Range("D44").Value = "1.834.54"
$ws.Range("E44").Value = "  +2.98%  "
$ws.Range("B45").Value = "BabyDogeCoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D45").Value = "0.0₈111"
$ws.Range("E45").Value = "  +3.93%  "
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").Value = "57.19"
$ws.Range("E46").Value = "  +5.68%  "
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "8.128"
$ws.Range("E47").Value = "  +2.94%  "
$ws.Range("B48").Value = "Frax"
$ws.Range("C48").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D48").Value = "1.002"
$ws.Range("E48").Value = "  -0.16%  "
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").Value = "0.05207"
$ws.Range("E49").Value = "  +0.18%  "
$ws.Range("B50").Value = "Aptos"
$ws.Range("C50").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D50").Value = "6.061"
$ws.Range("E50").Value = "  +4.47%  "
$ws.Range("B51").Value = "Mantle"
$ws.Range("C51").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D51").Value = "0.4243"
$ws.Range("E51").Value = "  +0.35%  "
